$wb = $excel.ActiveWorkbook

# Delete the car ("汽車") property sheet in its entirety.
$wb.Worksheets.Item("汽車").Delete()

# The source generator stamps every data row with a running "index" id
# (column A, and - on the Stock sheet only - duplicated again in column N).
# With the car rows gone, every id from the following sheets shifts down,
# so update each one in place to match the regenerated report.

$ws = $wb.Worksheets.Item("存款")
$ws.Range("A2").Value = 197
$ws.Range("A3").Value = 198
$ws.Range("A4").Value = 199
$ws.Range("A5").Value = 200
$ws.Range("A6").Value = 201
$ws.Range("A7").Value = 202
$ws.Range("A8").Value = 203
$ws.Range("A9").Value = 204
$ws.Range("A10").Value = 205
$ws.Range("A11").Value = 206
$ws.Range("A12").Value = 207
$ws.Range("A13").Value = 208
$ws.Range("A14").Value = 209
$ws.Range("A15").Value = 210
$ws.Range("A16").Value = 211
$ws.Range("A17").Value = 212
$ws.Range("A18").Value = 213
$ws.Range("A19").Value = 215
$ws.Range("A20").Value = 216
$ws.Range("A21").Value = 217
$ws.Range("A22").Value = 218
$ws.Range("A23").Value = 219
$ws.Range("A24").Value = 220
$ws.Range("A25").Value = 221
$ws.Range("A26").Value = 222
$ws.Range("A27").Value = 223
$ws.Range("A28").Value = 224
$ws.Range("A29").Value = 225
$ws.Range("A30").Value = 226
$ws.Range("A31").Value = 227
$ws.Range("A32").Value = 228
$ws.Range("A33").Value = 229
$ws.Range("A34").Value = 230
$ws.Range("A35").Value = 231

$ws = $wb.Worksheets.Item("股票")
$ws.Range("A2").Value = 238
$ws.Range("N2").Value = 238
$ws.Range("A3").Value = 239
$ws.Range("N3").Value = 239
$ws.Range("A4").Value = 240
$ws.Range("N4").Value = 240
$ws.Range("A5").Value = 241
$ws.Range("N5").Value = 241
$ws.Range("A6").Value = 242
$ws.Range("N6").Value = 242
$ws.Range("A7").Value = 243
$ws.Range("N7").Value = 243
$ws.Range("A8").Value = 244
$ws.Range("N8").Value = 244
$ws.Range("A9").Value = 246
$ws.Range("N9").Value = 246
$ws.Range("A10").Value = 247
$ws.Range("N10").Value = 247
$ws.Range("A11").Value = 248
$ws.Range("N11").Value = 248
$ws.Range("A12").Value = 249
$ws.Range("N12").Value = 249
$ws.Range("A13").Value = 250
$ws.Range("N13").Value = 250
$ws.Range("A14").Value = 251
$ws.Range("N14").Value = 251
$ws.Range("A15").Value = 252
$ws.Range("N15").Value = 252
$ws.Range("A16").Value = 253
$ws.Range("N16").Value = 253
$ws.Range("A17").Value = 254
$ws.Range("N17").Value = 254
$ws.Range("A18").Value = 255
$ws.Range("N18").Value = 255
$ws.Range("A19").Value = 256
$ws.Range("N19").Value = 256
$ws.Range("A20").Value = 257
$ws.Range("N20").Value = 257
$ws.Range("A21").Value = 258
$ws.Range("N21").Value = 258
$ws.Range("A22").Value = 259
$ws.Range("N22").Value = 259
$ws.Range("A23").Value = 260
$ws.Range("N23").Value = 260
$ws.Range("A24").Value = 261
$ws.Range("N24").Value = 261
$ws.Range("A25").Value = 262
$ws.Range("N25").Value = 262
$ws.Range("A26").Value = 263
$ws.Range("N26").Value = 263

$ws = $wb.Worksheets.Item("具有相當價值之財產")
$ws.Range("A2").Value = 290
$ws.Range("A3").Value = 291
$ws.Range("A4").Value = 292
$ws.Range("A5").Value = 293
$ws.Range("A6").Value = 294
$ws.Range("A7").Value = 295
$ws.Range("A8").Value = 296
$ws.Range("A9").Value = 297
$ws.Range("A10").Value = 298
$ws.Range("A11").Value = 299
$ws.Range("A12").Value = 300
$ws.Range("A13").Value = 301
$ws.Range("A14").Value = 302
$ws.Range("A15").Value = 303
$ws.Range("A16").Value = 304
$ws.Range("A17").Value = 305
$ws.Range("A18").Value = 307

$ws = $wb.Worksheets.Item("保險")
$ws.Range("A2").Value = 312
$ws.Range("A3").Value = 313
$ws.Range("A4").Value = 314
$ws.Range("A5").Value = 315
$ws.Range("A6").Value = 317
$ws.Range("A7").Value = 318
$ws.Range("A8").Value = 319

$ws = $wb.Worksheets.Item("債權")
$ws.Range("A2").Value = 324
$ws.Range("A3").Value = 325
$ws.Range("A4").Value = 326
$ws.Range("A5").Value = 327
$ws.Range("A6").Value = 328
$ws.Range("A7").Value = 329

$ws = $wb.Worksheets.Item("事業投資")
$ws.Range("A2").Value = 339
$ws.Range("A3").Value = 340
